$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A72").Value = "AOYH13"
$ws.Range("B72").Value = "Grasa para fusor"
$ws.Range("C72").Value = "15gr"
$ws.Range("D72").Value = 0
$ws.Range("E72").Value = 75000
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 5
$ws.Range("H72").Formula = "=(E72-D72)*G72"
$ws.Range("I72").Formula = "=D72*F72"
$ws.Range("J72").Value = 0
